$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: value changes; columns E and F are cleared (deleted entirely)
$rows = @(2, 3)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = 0.00497
    $ws.Range("E$r").ClearContents()
    $ws.Range("F$r").ClearContents()

    $ws.Range("G$r").Value = 0.3646967340590979
    $ws.Range("H$r").Value = 0.3646967340590979
    $ws.Range("I$r").Value = 0.1977708657335407
    $ws.Range("J$r").Value = 0.1724365812800706
    $ws.Range("K$r").Value = 63.3
    $ws.Range("L$r").Value = 0.1640746500777605
    $ws.Range("M$r").Value = 0.001
    $ws.Range("N$r").Value = 0.000003109452736318408
    $ws.Range("O$r").Value = 0.00001579778830963665
    $ws.Range("P$r").Value = 0.001
    $ws.Range("Q$r").Value = 0.000003109452736318408
    $ws.Range("R$r").Value = 0.00001579778830963665

    $ws.Range("U$r").Value = 143.6
    $ws.Range("V$r").Value = 0.4465174129353233
    $ws.Range("W$r").Value = 0.1785109983079526
    $ws.Range("X$r").Value = 0.05217253380938464
    $ws.Range("Y$r").Value = 0.126338464498568
    $ws.Range("Z$r").Value = 1.282451883123359
    $ws.Range("AA$r").Value = 0.2211416183819807
    $ws.Range("AB$r").Value = 0.04747123618731727
    $ws.Range("AC$r").Value = 0.1736703821946634
    $ws.Range("AD$r").Value = 62.7
    $ws.Range("AE$r").Value = 0
    $ws.Range("AF$r").Value = 62.7
    $ws.Range("AG$r").Value = -80.89999999999999
    $ws.Range("AH$r").Value = 0.1631537861046058
    $ws.Range("AI$r").Value = 0.1321948134092347
    $ws.Range("AJ$r").Value = -0.3361030328209388
    $ws.Range("AK$r").Value = -0.2446325975204112
    $ws.Range("AL$r").Value = 2.89
    $ws.Range("AM$r").Value = 2.89
    $ws.Range("AN$r").Value = 0.7116912599318956
    $ws.Range("AO$r").Value = 26.40138408304498
    $ws.Range("AP$r").Value = -0.9182746878547106
    $ws.Range("AQ$r").Value = 26.40138408304498
}
